$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the "Fecha / Volumen / Precio mínimo / Precio máximo /
# Precio promedio ponderado / Origen / Precio $/Kg" data among several rows
# (rows 4 and 12 are untouched). Target values per row, taken from the diff:
#   row -> (Fecha, Volumen, Precio, Origen)
$rows = @{
    2  = @{ D = 44592; M = 5;   Price = 7500; R = "Región de La Araucanía" }
    3  = @{ D = 44175; M = 40;  Price = 5000; R = "Provincia de Curicó" }
    5  = @{ D = 44999; M = 25;  Price = 2500; R = "Región de La Araucanía" }
    6  = @{ D = 44215; M = 65;  Price = 2800; R = "Región de La Araucanía" }
    7  = @{ D = 44214; M = 50;  Price = 1800; R = "Región de La Araucanía" }
    8  = @{ D = 44551; M = 120; Price = 4500; R = "Región de O'Higgins" }
    9  = @{ D = 44616; M = 200; Price = 3200; R = "Región de La Araucanía" }
    10 = @{ D = 44176; M = 20;  Price = 3000; R = "Región de O'Higgins" }
    11 = @{ D = 44998; M = 20;  Price = 2500; R = "Región de La Araucanía" }
    13 = @{ D = 44574; M = 200; Price = 3000; R = "Región de La Araucanía" }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.Price
    $ws.Range("O$r").Value = $vals.Price
    $ws.Range("P$r").Value = $vals.Price
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.Price
}
